$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H76").Value = 3339.1304
$ws.Range("I76").Value = 3309.0908
$ws.Range("K76").Value = 3309.0908
$ws.Range("M76").Value = -2994.0908
$ws.Range("H79").Value = 3339.1304
$ws.Range("I79").Value = 3309.0908
$ws.Range("K79").Value = 3309.0908
$ws.Range("M79").Value = -2217.0908
$ws.Range("H116").Value = 2803.5417
$ws.Range("I116").Value = 2108.75
$ws.Range("J116").Value = 3498.3333
$ws.Range("K116").Value = 2108.75
$ws.Range("L116").Value = 3498.3333
$ws.Range("M116").Value = 1333.25
$ws.Range("N116").Value = -10382.3333

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 8929.666999999999
$ws.Range("I32").Value = 10681.692
$ws.Range("K32").Value = 10681.692
$ws.Range("M32").Value = -10394.692
$ws.Range("H61").Value = 2610
$ws.Range("I61").Value = 2127.875
$ws.Range("J61").Value = 4263
$ws.Range("K61").Value = 2127.875
$ws.Range("L61").Value = 4263
$ws.Range("M61").Value = -1915.875
$ws.Range("N61").Value = -4687
$ws.Range("H63").Value = 2000
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 2000
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H74").Value = 3623.2327
$ws.Range("I74").Value = 829.95
$ws.Range("J74").Value = 6052.174
$ws.Range("K74").Value = 829.95
$ws.Range("L74").Value = 6052.174
$ws.Range("M74").Value = 44.04999999999995
$ws.Range("N74").Value = -7800.174
$ws.Range("H77").Value = 3623.2327
$ws.Range("I77").Value = 829.95
$ws.Range("J77").Value = 6052.174
$ws.Range("K77").Value = 4149.75
$ws.Range("L77").Value = 30260.87
$ws.Range("M77").Value = 218.25
$ws.Range("N77").Value = -38996.87
$ws.Range("H136").Value = 2610
$ws.Range("I136").Value = 2127.875
$ws.Range("J136").Value = 4263
$ws.Range("K136").Value = 6383.625
$ws.Range("L136").Value = 12789
$ws.Range("M136").Value = -3833.625
$ws.Range("N136").Value = -17889

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H86").Value = 1789.7333
$ws.Range("I86").Value = 1599.7037
$ws.Range("K86").Value = 1599.7037
$ws.Range("M86").Value = -476.7037
$ws.Range("H89").Value = 1789.7333
$ws.Range("I89").Value = 1599.7037
$ws.Range("K89").Value = 7998.5185
$ws.Range("M89").Value = -2382.5185

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H5").Value = 22728166
$ws.Range("I5").Value = 877.5454999999999
$ws.Range("J5").Value = 45455456
$ws.Range("K5").Value = 2632.6365
$ws.Range("L5").Value = 136366368
$ws.Range("M5").Value = -2520.6365
$ws.Range("N5").Value = -136366592
$ws.Range("H107").Value = 77346.80499999999
$ws.Range("I107").Value = 200310.6
$ws.Range("J107").Value = 48069.715
$ws.Range("K107").Value = 600931.8
$ws.Range("L107").Value = 144209.145
$ws.Range("M107").Value = -599011.8
$ws.Range("N107").Value = -148049.145
$ws.Range("H131").Value = 1925788.9
$ws.Range("I131").Value = 6596.125
$ws.Range("J131").Value = 2274733
$ws.Range("K131").Value = 19788.375
$ws.Range("L131").Value = 6824199
$ws.Range("M131").Value = -14748.375
$ws.Range("N131").Value = -6834279
$ws.Range("H135").Value = 22728166
$ws.Range("I135").Value = 877.5454999999999
$ws.Range("J135").Value = 45455456
$ws.Range("K135").Value = 7897.9095
$ws.Range("L135").Value = 409099104
$ws.Range("M135").Value = -5362.9095
$ws.Range("N135").Value = -409104174

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H80").Value = 88645.78999999999
$ws.Range("I80").Value = 1150
$ws.Range("J80").Value = 95376.234
$ws.Range("K80").Value = 1150
$ws.Range("L80").Value = 95376.234
$ws.Range("M80").Value = -152
$ws.Range("N80").Value = -97372.234
$ws.Range("H83").Value = 88645.78999999999
$ws.Range("I83").Value = 1150
$ws.Range("J83").Value = 95376.234
$ws.Range("K83").Value = 5750
$ws.Range("L83").Value = 476881.17
$ws.Range("M83").Value = -758
$ws.Range("N83").Value = -486865.17
$ws.Range("H122").Value = 2460.2
$ws.Range("I122").Value = 2559.4167
$ws.Range("J122").Value = 2063.3333
$ws.Range("K122").Value = 7678.250100000001
$ws.Range("L122").Value = 6189.999899999999
$ws.Range("M122").Value = -5228.250100000001
$ws.Range("N122").Value = -11089.9999
$ws.Range("H126").Value = 1742.7142
$ws.Range("I126").Value = 1279.8
$ws.Range("J126").Value = 2900
$ws.Range("K126").Value = 3839.4
$ws.Range("L126").Value = 8700
$ws.Range("M126").Value = -1369.4
$ws.Range("N126").Value = -13640

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H55").Value = 257.58334
$ws.Range("I55").Value = 240.42857
$ws.Range("J55").Value = 281.6
$ws.Range("K55").Value = 240.42857
$ws.Range("L55").Value = 281.6
$ws.Range("M55").Value = -67.42857000000001
$ws.Range("N55").Value = -627.6
$ws.Range("H82").Value = 1402.174
$ws.Range("I82").Value = 1157.6923
$ws.Range("J82").Value = 1720
$ws.Range("K82").Value = 1157.6923
$ws.Range("L82").Value = 1720
$ws.Range("M82").Value = -796.6922999999999
$ws.Range("N82").Value = -2442
$ws.Range("H85").Value = 1402.174
$ws.Range("I85").Value = 1157.6923
$ws.Range("J85").Value = 1720
$ws.Range("K85").Value = 1157.6923
$ws.Range("L85").Value = 1720
$ws.Range("M85").Value = 90.30770000000007
$ws.Range("N85").Value = -4216

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H62").Value = 12655.714
$ws.Range("I62").Value = 17247.5
$ws.Range("J62").Value = 6533.3335
$ws.Range("K62").Value = 17247.5
$ws.Range("L62").Value = 6533.3335
$ws.Range("M62").Value = -16623.5
$ws.Range("N62").Value = -7781.3335
$ws.Range("H65").Value = 12655.714
$ws.Range("I65").Value = 17247.5
$ws.Range("J65").Value = 6533.3335
$ws.Range("K65").Value = 86237.5
$ws.Range("L65").Value = 32666.6675
$ws.Range("M65").Value = -83117.5
$ws.Range("N65").Value = -38906.6675
$ws.Range("H96").Value = 1043.8334
$ws.Range("I96").Value = 1052.6
$ws.Range("J96").Value = 1000
$ws.Range("K96").Value = 1052.6
$ws.Range("L96").Value = 1000
$ws.Range("M96").Value = 320.4000000000001
$ws.Range("N96").Value = -3746
$ws.Range("H107").Value = 752.1539
$ws.Range("J107").Value = 908.6667
$ws.Range("L107").Value = 2726.0001
$ws.Range("N107").Value = -6566.0001
$ws.Range("H136").Value = 5253.1113
$ws.Range("I136").Value = 6789.5757
$ws.Range("J136").Value = 1027.8334
$ws.Range("K136").Value = 20368.7271
$ws.Range("L136").Value = 3083.5002
$ws.Range("M136").Value = -17818.7271
$ws.Range("N136").Value = -8183.5002
